# Adding Suppress Successful Logging and wbTypes in Config
# Applies the "Constants" sheet changes: removes the LogMessage_* rows and
# introduces the wb*_Type / wb*_SuppressSuccessful rows used by the
# REFramework logging workblocks.

$wb = $excel.ActiveWorkbook
$settings  = $wb.Worksheets.Item(1)   # Settings
$constants = $wb.Worksheets.Item(2)   # Constants
$assets    = $wb.Worksheets.Item(3)   # Assets

# Donor cell that already carries the "plain label" style (font #2 - plain
# Calibri 11 / default colour) so we can clone just the formatting.
$styleDonor = $settings.Range("B9")

function Set-WbRow($Row, $Name, $Value, $Description, $LeftAlign) {
    $constants.Cells.Item($Row, 1).Value = $Name
    $constants.Cells.Item($Row, 2).Value = $Value
    $constants.Cells.Item($Row, 3).Value = $Description

    if ($LeftAlign) {
        $valueCell = $constants.Cells.Item($Row, 2)
        $styleDonor.Copy()
        $valueCell.PasteSpecial(-4122) # xlPasteFormats
        $valueCell.HorizontalAlignment = -4131 # xlLeft
    }
}

Set-WbRow 19 "wbInit_Type" "MainTask, Initialization State" "Name of Workblock" $false
Set-WbRow 20 "wbInit_SuppressSuccessful" $false "Do not log successful executions of wb" $true
Set-WbRow 21 "wbGetTransactionData_Type" "MainTask, Get Transaction Data State" "Name of Workblock" $false
Set-WbRow 22 "wbGetTransactionData_SuppressSuccessful" $false "Do not log successful executions of wb" $true
Set-WbRow 23 "wbProcessTransaction_Type" "MainTask, Process Transaction State" "Name of Workblock" $false
Set-WbRow 24 "wbProcessTransaction_SuppressSuccessful" $false "Do not log successful executions of wb" $true
Set-WbRow 25 "wbNextTransaction_Type" "MainTask, Next Transaction" "Name of Workblock" $true
Set-WbRow 26 "wbNextTransaction_SuppressSuccessful" $false "Do not log successful executions of wb" $true
Set-WbRow 27 "wbCloseAllApplications_Type" "MainTask, Close All Applications" "Name of Workblock" $false
Set-WbRow 28 "wbCloseAllApplications_SuppressSuccessful" $false "Do not log successful executions of wb" $true
Set-WbRow 29 "wbInitAllApplications_Type" "MainTask, InitAllApplications" "Name of Workblock" $false
Set-WbRow 30 "wbInitAllApplications_SuppressSuccessful" $false "Do not log successful executions of wb" $true
Set-WbRow 31 "wbProcess_Type" "MainTask, Process Transaction" "Name of Workblock" $false
Set-WbRow 32 "wbProcess_SuppressSuccessful" $false "Do not log successful executions of wb" $true

# The Constants tab is the active tab, with C24 selected, matching the
# authored workbook state.
$constants.Activate()
$constants.Range("C24").Select()
